$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.158.15'
$ws.Range('E2').Value = '  +8.66%  '
$ws.Range('D3').Value = '3.523.45'
$ws.Range('E3').Value = '  +12.69%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '190.41'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +13.96%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '552.50'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +8.21%  '
$ws.Range('D7').Value = '3.516.76'
$ws.Range('E7').Value = '  +12.46%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.607'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.41%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.634'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +8.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.150'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +19.65%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '55.24'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.75%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000269'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +10.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.40'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +7.48%  '
$ws.Range('D15').Value = '4.063.52'
$ws.Range('E15').Value = '  +12.67%  '
$ws.Range('D16').Value = '3.517.23'
$ws.Range('E16').Value = '  +13.15%  '
$ws.Range('E17').Value = '  +6.99%  '
$ws.Range('D18').Value = '67.222.59'
$ws.Range('E18').Value = '  +9.37%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.21'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +9.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.87'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +11.10%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.996'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +6.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '429.84'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +20.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.93'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.74%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.29'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +8.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.22'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +10.99%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.20'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +15.14%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.17'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.98'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +10.13%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.90'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +12.64%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '30.28'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +9.93%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '654.60'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.96%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.68'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.73'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +6.46%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.111'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +9.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '59.61'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +6.64%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '38.81'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +10.48%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0818'
$ws.Range('E38').Value = '  +22.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.392'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.08%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.139'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +15.79%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.30'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +15.20%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '3.023.39'
$ws.Range('E44').Value = '  +7.18%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.67'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +5.75%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.89'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +18.05%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.29'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +11.89%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0418'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +10.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.72'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +4.40%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.131'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +9.19%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.79'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +19.85%  '
